$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8076208178438662
$ws.Range("B3").Value = 0.7358490566037735
$ws.Range("B4").Value = 0.8154639175257732
$ws.Range("B5").Value = 0.3035019455252918
$ws.Range("B6").Value = 0.9658119658119658
$ws.Range("B7").Value = 0.8595555339428127

$wb.Save()
